$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Magic Player Rewards 2008 (P08)",
    "Desert",
    "Eternal Witness",
    "Isochron Scepter",
    "Pendelhaven",
    "Remand",
    "Resurrection",
    "Serrated Arrows",
    "Shrapnel Blast",
    "Tendrils of Agony",
    "Thirst for Knowledge",
    "Tormod's Crypt",
    "Wall of Roots"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
